$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-group permutations: each row in a group receives the B:AC content
# that another row in the same group held before the edit (ids/odds were
# re-paired while keeping the row index/order column A and date column E).

# Group 1: rows 23, 24
$row23 = $ws.Range("B23:AC23").Value2
$row24 = $ws.Range("B24:AC24").Value2
$ws.Range("B23:AC23").Value2 = $row24
$ws.Range("B24:AC24").Value2 = $row23

# Group 2: rows 38, 39
$row38 = $ws.Range("B38:AC38").Value2
$row39 = $ws.Range("B39:AC39").Value2
$ws.Range("B38:AC38").Value2 = $row39
$ws.Range("B39:AC39").Value2 = $row38

# Group 3: rows 49, 50
$row49 = $ws.Range("B49:AC49").Value2
$row50 = $ws.Range("B50:AC50").Value2
$ws.Range("B49:AC49").Value2 = $row50
$ws.Range("B50:AC50").Value2 = $row49

# Group 4: rows 59, 61, 60
$row59 = $ws.Range("B59:AC59").Value2
$row61 = $ws.Range("B61:AC61").Value2
$row60 = $ws.Range("B60:AC60").Value2
$ws.Range("B59:AC59").Value2 = $row61
$ws.Range("B61:AC61").Value2 = $row60
$ws.Range("B60:AC60").Value2 = $row59

# Group 5: rows 67, 68
$row67 = $ws.Range("B67:AC67").Value2
$row68 = $ws.Range("B68:AC68").Value2
$ws.Range("B67:AC67").Value2 = $row68
$ws.Range("B68:AC68").Value2 = $row67

# Group 6: rows 79, 80
$row79 = $ws.Range("B79:AC79").Value2
$row80 = $ws.Range("B80:AC80").Value2
$ws.Range("B79:AC79").Value2 = $row80
$ws.Range("B80:AC80").Value2 = $row79

# Group 7: rows 97, 99
$row97 = $ws.Range("B97:AC97").Value2
$row99 = $ws.Range("B99:AC99").Value2
$ws.Range("B97:AC97").Value2 = $row99
$ws.Range("B99:AC99").Value2 = $row97

# Group 8: rows 102, 103
$row102 = $ws.Range("B102:AC102").Value2
$row103 = $ws.Range("B103:AC103").Value2
$ws.Range("B102:AC102").Value2 = $row103
$ws.Range("B103:AC103").Value2 = $row102

# Group 9: rows 107, 109
$row107 = $ws.Range("B107:AC107").Value2
$row109 = $ws.Range("B109:AC109").Value2
$ws.Range("B107:AC107").Value2 = $row109
$ws.Range("B109:AC109").Value2 = $row107

# Group 10: rows 114, 115
$row114 = $ws.Range("B114:AC114").Value2
$row115 = $ws.Range("B115:AC115").Value2
$ws.Range("B114:AC114").Value2 = $row115
$ws.Range("B115:AC115").Value2 = $row114

# Group 11: rows 121, 122
$row121 = $ws.Range("B121:AC121").Value2
$row122 = $ws.Range("B122:AC122").Value2
$ws.Range("B121:AC121").Value2 = $row122
$ws.Range("B122:AC122").Value2 = $row121

# Group 12: rows 125, 127, 126
$row125 = $ws.Range("B125:AC125").Value2
$row127 = $ws.Range("B127:AC127").Value2
$row126 = $ws.Range("B126:AC126").Value2
$ws.Range("B125:AC125").Value2 = $row127
$ws.Range("B127:AC127").Value2 = $row126
$ws.Range("B126:AC126").Value2 = $row125

# Group 13: rows 145, 146
$row145 = $ws.Range("B145:AC145").Value2
$row146 = $ws.Range("B146:AC146").Value2
$ws.Range("B145:AC145").Value2 = $row146
$ws.Range("B146:AC146").Value2 = $row145

# Group 14: rows 150, 151
$row150 = $ws.Range("B150:AC150").Value2
$row151 = $ws.Range("B151:AC151").Value2
$ws.Range("B150:AC150").Value2 = $row151
$ws.Range("B151:AC151").Value2 = $row150

# Group 15: rows 156, 157
$row156 = $ws.Range("B156:AC156").Value2
$row157 = $ws.Range("B157:AC157").Value2
$ws.Range("B156:AC156").Value2 = $row157
$ws.Range("B157:AC157").Value2 = $row156

# Group 16: rows 184, 185
$row184 = $ws.Range("B184:AC184").Value2
$row185 = $ws.Range("B185:AC185").Value2
$ws.Range("B184:AC184").Value2 = $row185
$ws.Range("B185:AC185").Value2 = $row184

# Group 17: rows 192, 193, 194
$row192 = $ws.Range("B192:AC192").Value2
$row193 = $ws.Range("B193:AC193").Value2
$row194 = $ws.Range("B194:AC194").Value2
$ws.Range("B192:AC192").Value2 = $row193
$ws.Range("B193:AC193").Value2 = $row194
$ws.Range("B194:AC194").Value2 = $row192
